# Fixing Bulk Operation Template: insert a "Group" column between
# "Category" and "Unit" on the "Create Item" sheet, and set the sheet
# to print in portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D (shifts old Unit/Sell Price from D/E to E/F)
# and give it the same header formatting/width as its left neighbour
# (Category, column C).
$ws.Columns("D:D").Insert(-4161)

$ws.Range("D1").Value = "Group"

# Match the other header cells' width/style so B:D end up uniform.
$ws.Range("B1:D1").ColumnWidth = $ws.Range("C1").ColumnWidth

# Print in portrait orientation.
$ws.PageSetup.Orientation = 1
